# Natmi following Dr Hou advice
#
# The original sheet only had 3 cell-communication rows (FAPs->Pth1r->{ECs,FAPs,sCs}).
# The updated analysis adds a new "M2" target cluster and a new "sCs" sending
# cluster, expanding the Pthlh/Pth1r edge table from 3 rows to 8 rows
# (sheet rows 2-9), each with refreshed expression/specificity statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster | B Ligand symbol | C Receptor symbol | D Target cluster
# E Ligand-expressing cells | F Ligand detection rate
# G Ligand average expression value | H Ligand total expression value
# I Ligand derived specificity of average expression value
# J Ligand derived specificity of total expression value
# K Receptor-expressing cells | L Receptor detection rate
# M Receptor average expression value | N Receptor total expression value
# O Receptor derived specificity of average expression value
# P Receptor derived specificity of total expression value
# Q Edge average expression weight | R Edge total expression weight
# S Edge average expression derived specificity | T Edge total expression derived specificity

$data = @(
        @("FAPs", "Pthlh", "Pth1r", "ECs",  3, 1, 6.802404,            20.407212,           0.9280918435277148,   0.9280918435277147,   3, 1, 1.240767,            3.722301,   0.1653663592679705,   0.1653663592679705,   8.440198403867999,  75.961785634812,    0.1534751692304771,    0.1534751692304771),
        @("FAPs", "Pthlh", "Pth1r", "FAPs", 3, 1, 6.802404,            20.407212,           0.9280918435277148,   0.9280918435277147,   3, 1, 4.697042333333333,  14.091127,  0.6260101936873453,   0.6260101936873453,   31.951179556436,    287.560616007924,   0.5809949547264301,    0.58099495472643),
        @("FAPs", "Pthlh", "Pth1r", "M2",   3, 1, 6.802404,            20.407212,           0.9280918435277148,   0.9280918435277147,   1, 0.3333333333333333, 0.072744,           0.218232,   0.009695140536933402, 0.009695140536933402, 0.494834076576,     4.453506689184,     0.0089979808541828,     0.008997980854182798),
        @("FAPs", "Pthlh", "Pth1r", "sCs",  3, 1, 6.802404,            20.407212,           0.9280918435277148,   0.9280918435277147,   3, 1, 1.492587,            4.477761,   0.1989283065077507,   0.1989283065077507,   10.153179779148,    91.37861801233201,  0.1846237387166247,    0.1846237387166246),
        @("sCs",  "Pthlh", "Pth1r", "ECs",  1, 0.3333333333333333, 0.5270473333333333, 1.581142,            0.07190815647228528,  0.07190815647228528,  3, 1, 1.240767,            3.722301,   0.1653663592679705,   0.1653663592679705,   0.6539429386379999, 5.885486447742,     0.01189119003749337,   0.01189119003749336),
        @("sCs",  "Pthlh", "Pth1r", "FAPs", 1, 0.3333333333333333, 0.5270473333333333, 1.581142,            0.07190815647228528,  0.07190815647228528,  3, 1, 4.697042333333333,  14.091127,  0.6260101936873453,   0.6260101936873453,   2.475563636337111,  22.280072727034,    0.04501523896091524,   0.04501523896091524),
        @("sCs",  "Pthlh", "Pth1r", "M2",   1, 0.3333333333333333, 0.5270473333333333, 1.581142,            0.07190815647228528,  0.07190815647228528,  1, 0.3333333333333333, 0.072744,           0.218232,   0.009695140536933402, 0.009695140536933402, 0.038339531216,     0.345055780944,     0.0006971596827506031,  0.0006971596827506031),
        @("sCs",  "Pthlh", "Pth1r", "sCs",  1, 0.3333333333333333, 0.5270473333333333, 1.581142,            0.07190815647228528,  0.07190815647228528,  3, 1, 1.492587,            4.477761,   0.1989283065077507,   0.1989283065077507,   0.786663998118,     7.079975983062,     0.01430456779112606,   0.01430456779112606)
)

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $columns.Length; $j++) {
        $ws.Range($columns[$j] + "$rowNum").Value = $rowValues[$j]
    }
}

Write-Output "updated Pthlh-Pth1r sheet to $($data.Length) data rows"
